$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 523326
$ws.Range("R2").Value = 6619778
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
